$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.104.32"
$ws.Range("E2").Value = "  +2.35%  "

$ws.Range("D3").Value = "2.302.83"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.56"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.25"
$ws.Range("E6").Value = "  +5.21%  "

$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +3.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.28"
$ws.Range("E10").Value = "  +4.29%  "

$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.12"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("E13").Value = "  +4.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.86"
$ws.Range("E14").Value = "  +16.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.81"
$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("D16").Value = "2.662.60"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").Value = "2.335.47"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.810"
$ws.Range("E18").Value = "  +4.56%  "

$ws.Range("D19").Value = "42.981.31"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("E20").Value = "  +8.81%  "

$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  +2.31%  "

$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.85"
$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.28"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  +10.66%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.62"
$ws.Range("E28").Value = "  +3.55%  "

$ws.Range("E29").Value = "  +11.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.04"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.02"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.19"
$ws.Range("E32").Value = "  +1.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.03"
$ws.Range("E34").Value = "  +2.00%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.56"
$ws.Range("E35").Value = "  +4.99%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("E36").Value = "  +4.29%  "

$ws.Range("E37").Value = "  +7.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0699"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("E39").Value = "  +3.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"

$ws.Range("E41").Value = "  +5.14%  "

$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("D44").Value = "1.998.47"

$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.07"
$ws.Range("E46").Value = "  +5.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.68"
$ws.Range("E47").Value = "  +2.31%  "

$ws.Range("E48").Value = "  +3.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.67"
$ws.Range("E49").Value = "  +4.74%  "

$ws.Range("D50").Value = "2.528.22"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("E51").Value = "  +3.21%  "
